# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The "Periodo Mora" list (rows 16-54, column E) used to run descending
# from 2003 down to 1701. It is now reordered to run ascending from 1701
# up to 2003. Along with the reorder, the "Salario Basico" (F) /
# "Valor Mora" (G) figures for the block were refreshed: every row's
# Valor Mora becomes 781242 (was 689500), and the Salario Basico split
# flips - the first 20 periods (1701..1808) now carry 27580 while the
# remaining 19 (1809..2003) carry 31249 (previously it was the other way
# around).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$periods = @(
    "1701", "1702", "1703", "1704", "1705", "1706", "1707", "1708", "1709", "1710",
    "1711", "1712", "1801", "1802", "1803", "1804", "1805", "1806", "1807", "1808",
    "1809", "1810", "1811", "1812", "1901", "1902", "1903", "1904", "1905", "1906",
    "1907", "1908", "1909", "1910", "1911", "1912", "2001", "2002", "2003"
)

$firstRow = 16
$newValorMora = 781242

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $firstRow + $i

    if ($i -lt 20) {
        $salarioBasico = 27580
    } else {
        $salarioBasico = 31249
    }

    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $salarioBasico
    $ws.Range("G$row").Value = $newValorMora
}
